# Updates the cryptos list (Price / Volume(1h) columns) with a fresh
# scrape, as run by the GitHub Actions workflow on
# Wed Nov 29 11:43:50 UTC 2023. Also reflects Aave and InjectiveProtocol
# swapping rank positions (rows 41/42) in the refreshed data.
#
# Price cells that look like a bare decimal number (e.g. "229.35") are
# pinned to Text format first so Excel's automatic type detection doesn't
# turn them into floating point numbers -- the source data stores prices
# as plain strings (note some entries use "." as a thousands separator,
# e.g. "38.163.78", which would not parse as a number anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.163.78"
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").Value = "2.056.13"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.35"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.96"
$ws.Range("E7").Value = "  +8.66%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.79"
$ws.Range("E12").Value = "  +3.49%  "
$ws.Range("D13").Value = "2.360.73"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("E14").Value = "  +5.03%  "
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D17").Value = "2.057.50"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "38.091.95"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.30"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.86"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.75"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.75"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.23"
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.133"
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.96"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.60"
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("E34").Value = "  +7.84%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.40"
$ws.Range("E36").Value = "  +16.62%  "
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").Value = "1.518.86"
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.66"
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.01"
$ws.Range("E42").Value = "  +4.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0926"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.05"
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").Value = "2.248.53"
$ws.Range("E51").Value = "  +1.63%  "
